$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Sprint numbers (column C, "Sprint #")
$ws.Range("C3").Value = 8
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 3
$ws.Range("C6").Value = 4
$ws.Range("C7").Value = 5
$ws.Range("C8").Value = 4
$ws.Range("C9").Value = 5
$ws.Range("C10").Value = 7
$ws.Range("C11").Value = 7
$ws.Range("C12").Value = 5
$ws.Range("C13").Value = 8
$ws.Range("C14").Value = 3
$ws.Range("C15").Value = 6
$ws.Range("C16").Value = 4
$ws.Range("C17").Value = 8
$ws.Range("C18").Value = 7
$ws.Range("C19").Value = 2
$ws.Range("C20").Value = 5
$ws.Range("C21").Value = 5
$ws.Range("C22").Value = 6
$ws.Range("C23").Value = 6
$ws.Range("C24").Value = 7
$ws.Range("C25").Value = 5
$ws.Range("C26").Value = 7
$ws.Range("C27").Value = 7
$ws.Range("C28").Value = 7
$ws.Range("C29").Value = 8

# Move the active selection to C29 as in the saved workbook view
$ws.Range("C29").Select()
